# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row swaps (coin ranking order changed) ---
$ws.Range('B9').Value = 'Solana'
$ws.Range('C9').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.00'
$ws.Range('E9').Value = '  +2.09%  '

$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2896'
$ws.Range('E10').Value = '  -0.29%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01842'
$ws.Range('E37').Value = '  +0.51%  '

$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.938'
$ws.Range('E38').Value = '  +2.28%  '

# --- Price (column D) updates ---
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.400.51'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.842.07'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.35'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6266'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07444'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07719'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.844.22'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.977'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001030'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.82'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.239'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.460.15'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '232.70'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.329'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.08'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.495'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1351'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.35'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.07139'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.465'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.484'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.045'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.035'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.140'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6990'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.577'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.817'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.235.11'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9615'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.0000'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.014.24'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.93'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.52'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.972'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.898'

# --- Volume(1h) (column E) updates ---
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('E15').Value = '  -2.27%  '
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('E19').Value = '  +1.51%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('E22').Value = '  -2.21%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('E26').Value = '  -1.76%  '
$ws.Range('E27').Value = '  -1.22%  '
$ws.Range('E28').Value = '  +11.18%  '
$ws.Range('E29').Value = '  +3.44%  '
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('E31').Value = '  -1.28%  '
$ws.Range('E32').Value = '  -1.52%  '
$ws.Range('E33').Value = '  -0.68%  '
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('E39').Value = '  -1.04%  '
$ws.Range('E40').Value = '  -2.75%  '
$ws.Range('E41').Value = '  +5.39%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('E45').Value = '  -1.30%  '
$ws.Range('E46').Value = '  +3.70%  '
$ws.Range('E47').Value = '  -0.70%  '
$ws.Range('E49').Value = '  -1.29%  '
$ws.Range('E50').Value = '  -2.62%  '
$ws.Range('E51').Value = '  -1.73%  '

Write-Host 'cryptos list updated'
